# Rescale EV power to fix error, ev allocation as tuples
#
# This script updates specific cells in the p_mw.xlsx time-series sheet
# that represent EV (electric vehicle) load values which were previously
# zero but should now carry a small non-zero power allocation (0.0066,
# with 0.0132 where two EV allocations overlap in the same cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("W9").Value = 0.0066

# Row 11
$ws.Range("V11").Value = 0.0066
$ws.Range("X11").Value = 0.0066

# Row 12
$ws.Range("X12").Value = 0.0066

# Row 13
$ws.Range("V13").Value = 0.0066

# Row 14
$ws.Range("V14").Value = 0.0066

# Row 15
$ws.Range("V15").Value = 0.0066

# Row 17
$ws.Range("W17").Value = 0.0066

# Row 19
$ws.Range("U19").Value = 0.0066
$ws.Range("W19").Value = 0.0132
$ws.Range("X19").Value = 0.0066

# Row 20
$ws.Range("W20").Value = 0.0066
$ws.Range("X20").Value = 0.0066

# Row 21
$ws.Range("V21").Value = 0.0066

# Row 22
$ws.Range("V22").Value = 0.0066
